$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.291.60"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "2.512.30"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.08"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.35"
$ws.Range("E6").Value = "  -1.86%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("E8").Value = "  +1.67%  "
$ws.Range("D9").Value = "2.510.86"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  -1.64%  "
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.357"
$ws.Range("E12").Value = "  +4.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.92"
$ws.Range("E13").Value = "  +2.68%  "
$ws.Range("D14").Value = "2.976.96"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "68.972.35"
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("E16").Value = "  -2.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.87"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "2.514.51"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.36"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.70"
$ws.Range("E20").Value = "  +2.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.28"
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.94"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.02"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("E26").Value = "  -1.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.93"
$ws.Range("E27").Value = "  -2.87%  "
$ws.Range("D28").Value = "2.650.40"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "464.09"
$ws.Range("E32").Value = "  -3.07%  "
$ws.Range("E33").Value = "  -2.00%  "
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.117"
$ws.Range("E36").Value = "  +0.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.58"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.01"
$ws.Range("E38").Value = "  +0.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.54"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  +1.17%  "
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("E43").Value = "  -2.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.30"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  -5.54%  "
$ws.Range("E46").Value = "  -13.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.70"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0729"
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("E51").Value = "  -3.39%  "
